$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.917.04'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +5.71%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.248.73'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +3.96%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '230.05'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.94%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.36%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '61.35'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -3.54%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.407'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '58.61'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +1.05%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0884'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +4.08%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.26%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.583.54'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +3.98%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.85'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.16%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.81'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.90%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.806'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.53%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +2.22%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.248.55'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +4.03%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '41.787.32'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +5.55%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +2.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.14'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.31%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0896'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +5.93%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '249.98'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +9.62%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.07%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.39'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.73%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.26'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -4.70%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.63'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.31%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '168.09'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.52%  '
$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.142'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +2.18%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +2.20%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.44'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +2.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.79'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +3.86%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.69%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.20'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +10.90%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.74'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +2.92%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.55%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.84'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +5.49%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -3.41%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.40%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.000251'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +33.91%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.03'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +3.30%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0237'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +4.46%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +12.51%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '100.51'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.59%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0980'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +6.14%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.486.27'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.79%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -3.11%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '16.55'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -6.14%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.31%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.78'
